$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header + data columns (E:G) for the Profile Picture Upload scenario ---

# Copy header style (bold, bordered, centered - same as D1) into E1:G1
$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("D1").Copy($ws.Range("F1"))
$ws.Range("D1").Copy($ws.Range("G1"))

# Copy data-row style (bordered, centered - same as D2) into E2:G2
$ws.Range("D2").Copy($ws.Range("E2"))
$ws.Range("D2").Copy($ws.Range("F2"))
$ws.Range("D2").Copy($ws.Range("G2"))

# Header text
$ws.Range("E1").Value = "Exceeded 1MB"
$ws.Range("F1").Value = "Incorrect Dimensions"
$ws.Range("G1").Value = "Unsupported File Format"

# Data values
$ws.Range("E2").Value = "D:\MARINA\SDET\GroupProject\Exceeded 1MB.jpg"
$ws.Range("F2").Value = "D:\MARINA\SDET\GroupProject\Incorrect Dimensions.jpg"
$ws.Range("G2").Value = "D:\MARINA\SDET\GroupProject\Unsupported FIle Format.bmp"

# --- Column widths ---
# Columns A:C narrow from 16.26953125 to 13 characters
$ws.Range("A1:C2").ColumnWidth = 12.166666666666666
# Column D widens from 49.26953125 to 52.90625 characters
$ws.Range("D1:D2").ColumnWidth = 52
# New columns E:G use the same display width as D (no explicit style index)
$ws.Range("E1:G2").ColumnWidth = 52

# --- Selection moves to D10 ---
$ws.Range("D10").Select() | Out-Null
